$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repull data, push all data, mean calculation:
# Update the "dSF" (column F) values for the affected rows.
$ws.Range("F3").Value = -3
$ws.Range("F13").Value = 1
$ws.Range("F15").Value = -1
$ws.Range("F20").Value = 2
$ws.Range("F24").Value = -3
$ws.Range("F27").Value = 1
$ws.Range("F37").Value = 1
$ws.Range("F40").Value = 0
$ws.Range("F44").Value = 0
$ws.Range("F55").Value = 4
$ws.Range("F62").Value = 7
$ws.Range("F64").Value = -7
$ws.Range("F68").Value = 7
